# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a handback event:
#   - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Handback DateTime values are refreshed
#   - Error Detail (stale-handback) messages are cleared now that files are in sync
#   - The Status / Error Detail columns are resized

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Refresh the handback timestamps
$zhcn.Range("K2").Value = "2016-08-12 12:40:58"
$zhcn.Range("K3").Value = "2016-08-12 12:40:58"

# The handback files are now in sync -> clear the stale-handback error detail
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# Resize Status / Error Detail columns
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Refresh the handback timestamps
$dede.Range("K2").Value = "2016-08-12 12:41:12"
$dede.Range("K3").Value = "2016-08-12 12:41:12"

# The handback files are now in sync -> clear the stale-handback error detail
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# Resize Status / Error Detail columns
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
